# Update on 2018-05-30, 支出生活费300
# Adds a new expense entry row (row 41) to the detail table:
#   序号(B)=39, 类别(C)=支出, 金额(D)=300, 时间(E)=2018-05-30,
#   费用类别(F)=生活费, 备注(G)=生活费(5/31-6/9)
# J3/K3/J9 are existing SUMIFS summary formulas that recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 was a blank placeholder row styled like the data rows above it
# (row 40). Copy that row's number/date/text formatting down onto row 41
# before filling in the new entry so the new cells keep the same look
# (fill, border, date number format, left-aligned text) as the rest of
# the table instead of falling back to the row's bare placeholder style.
$ws.Range("D40:G40").Copy()
$ws.Range("D41:G41").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B41").Value = 39
$ws.Range("C41").Value = "支出"
$ws.Range("D41").Value = 300
$ws.Range("E41").Value = 43250
$ws.Range("F41").Value = "生活费"
$ws.Range("G41").Value = "生活费(5/31-6/9)"

# Match the author's resulting scroll position / active cell on save.
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B41").Select()

$wb.Application.Calculate()
